# Regenerate the save_data "K" column (column G) values for this pitcher's
# game log. The K column previously held "Strike#"-style counts; it is
# recomputed here (std/mean regenerated upstream, s_vals recalculated) and
# the resulting values are written back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value, as produced by the regen step.
$kValues = [ordered]@{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 0
    7  = 2
    8  = 3
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 4
    16 = 3
    17 = 1
    18 = 3
    20 = 3
    21 = 1
    22 = 4
    23 = 0
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 1
    34 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
